$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("C2").Value = "南宁·第一届ANE·DACG动漫嘉年华（取消）"
$ws1.Range("F2").Value = 1038
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F3").Value = 2201
$ws1.Range("F5").Value = 480

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("C4").Value = "南宁·第一届ANE·DACG动漫嘉年华（取消）"
$ws4.Range("F4").Value = 1038
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F5").Value = 2201
$ws4.Range("F7").Value = 480
